$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-04 Monday" "2024-03-05 Tuesday"

Replace-Text "42÷7=6, 0" "13÷5=2, 3"
Replace-Text "22÷8=2, 6" "85÷3=28, 1"
Replace-Text "37÷3=12, 1" "13÷7=1, 6"
Replace-Text "51÷9=5, 6" "23÷2=11, 1"
Replace-Text "21÷6=3, 3" "27÷9=3, 0"
Replace-Text "75÷6=12, 3" "45÷6=7, 3"
Replace-Text "90÷3=30, 0" "49÷5=9, 4"
Replace-Text "97÷7=13, 6" "26÷9=2, 8"
Replace-Text "89÷5=17, 4" "11÷8=1, 3"
Replace-Text "22÷6=3, 4" "72÷7=10, 2"
Replace-Text "10÷4=2, 2" "70÷9=7, 7"
Replace-Text "68÷7=9, 5" "79÷4=19, 3"
Replace-Text "56÷2=28, 0" "37÷8=4, 5"
Replace-Text "93÷2=46, 1" "29÷7=4, 1"
Replace-Text "84÷4=21, 0" "88÷6=14, 4"
Replace-Text "30÷4=7, 2" "69÷6=11, 3"
Replace-Text "82÷8=10, 2" "75÷5=15, 0"
Replace-Text "50÷6=8, 2" "80÷9=8, 8"
Replace-Text "92÷7=13, 1" "19÷4=4, 3"
Replace-Text "62÷2=31, 0" "40÷4=10, 0"
Replace-Text "58÷5=11, 3" "50÷8=6, 2"
Replace-Text "65÷2=32, 1" "28÷3=9, 1"
Replace-Text "41÷7=5, 6" "43÷6=7, 1"
Replace-Text "39÷7=5, 4" "70÷6=11, 4"
Replace-Text "31÷5=6, 1" "53÷5=10, 3"

Write-Output "Done replacing text"
